$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.163.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.920.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3009"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06920"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.914.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07316"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6857"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.106"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.100.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008086"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.160.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9976"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.891"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "179.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +31.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.139"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.380"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.959"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.409"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.378"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08988"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05271"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7568"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.663"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01924"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.747"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9436"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.41%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4384"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.44%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.955"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.50%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.857"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.50%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1343"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05867"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.624"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.82%  "

$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3907"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.36%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.395"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.53%  "

